# Applies the "new .ttl from Google sheet has been generated" update to the
# vocabulary worksheet: refresh the base URI, clear out the example
# title/description/creator metadata and the example term rows, shift the
# header block up, and drop the now-unused trailing blank template rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Base URI changed from the m4m-dk-3 namespace to the test namespace ---
$ws.Range("B1").Value = "http://purl.org/test/variables/"
$ws.Range("C3").Value = "http://purl.org/test/variables/"

# --- 2. Clear the example title / description values (labels stay) ---
$ws.Range("B10").ClearContents()
$ws.Range("B11").ClearContents()

# --- 3. Clear the example creator rows (12-16), keep only the first dct:creator label ---
$ws.Range("B12:C12").ClearContents()
$ws.Range("B13:C16").ClearContents()

# --- 4. Rows 13-16 become the dct:rights / pav:version / pav:createdOn / pav:lastUpdatedOn rows ---
$ws.Range("A13").Value = "dct:rights"
$ws.Range("C13").Value = "License under which the vocabulary is provided"

$ws.Range("A14").Value = "pav:version"
$ws.Range("C14").Value = "Vocabulary version"

$ws.Range("A15").Value = "pav:createdOn"
$ws.Range("C15").Value = "Date when vocabulary was initially created (follow https://en.wikipedia.org/wiki/ISO_8601)"

$ws.Range("A16").Value = "pav:lastUpdatedOn"
$ws.Range("C16").Value = "Date of the last vocabulary update"

# --- 5. Row 17 becomes the "Definition of terms" section header ---
$ws.Range("A17").Value = "Definition of terms (optionally properties)"

# --- 6. Row 18 becomes the "Identifier" column-header row (moved up from row 23) ---
$ws.Range("A18").Value = "Identifier"
$ws.Range("B18").Value = "skos:prefLabel@en"
$ws.Range("C18").Value = "qudt:unit(separator=" + [char]34 + "," + [char]34 + ")"
$ws.Range("D18").Value = "skos:altLabel(separator=" + [char]34 + ";" + [char]34 + ")"
$ws.Range("E18").Value = "skos:definition@en"
$ws.Range("F18").Value = "dct:source(separator=" + [char]34 + "," + [char]34 + ")"
$ws.Range("G18").Value = "skos:broader(separator=" + [char]34 + "," + [char]34 + ")"
$ws.Range("H18").Value = "skos:exactMatch(separator=" + [char]34 + "," + [char]34 + ")"
$ws.Range("I18").Value = "skos:closeMatch(separator=" + [char]34 + "," + [char]34 + ")"
$ws.Range("J18").Value = "skos:editorialNote@en"
$ws.Range("K18").Value = "dct:creator(separator=" + [char]34 + "," + [char]34 + ")"
$ws.Range("L18").Value = "dct:contributor(separator=" + [char]34 + "," + [char]34 + ")"

# --- 7. Rows 19-32 become blank "vars:" template rows (old sample term data removed) ---
for ($r = 19; $r -le 32; $r++) {
    $ws.Range("A" + $r).Value = "vars:"
    $ws.Range("B" + $r + ":T" + $r).ClearContents()
}

# --- 8. Drop the now-superfluous trailing blank template rows 90-94 ---
$ws.Range("A90:T94").EntireRow.Delete()
